$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the crypto price refresh diff.
# Numeric-looking price strings need NumberFormat "@" forced first so Excel
# keeps them as literal text (preserving trailing zeros / multi-dot values)
# instead of auto-converting them to numbers.

$ws.Range('D2').Value = '90.708.15'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '3.113.10'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.30'
$ws.Range('E5').Value = '  +9.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '631.02'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  -1.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.360'
$ws.Range('E8').Value = '  -3.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = '3.111.90'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.713'
$ws.Range('E11').Value = '  -7.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.195'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.70'
$ws.Range('E13').Value = '  +5.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000244'
$ws.Range('E14').Value = '  -2.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.47'
$ws.Range('E15').Value = '  -2.07%  '
$ws.Range('D16').Value = '90.615.07'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').Value = '3.712.45'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').Value = '3.112.84'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.12'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000207'
$ws.Range('E21').Value = '  -4.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '441.06'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.54'
$ws.Range('E23').Value = '  +5.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.93'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.71'
$ws.Range('E25').Value = '  -9.95%  '
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.53'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '88.18'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').Value = '3.307.53'
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.55'
$ws.Range('E30').Value = '  +3.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.159'
$ws.Range('E31').Value = '  -3.78%  '
$ws.Range('E32').Value = '  +7.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.26'
$ws.Range('E33').Value = '  +4.11%  '
$ws.Range('E34').Value = '  +18.85%  '
$ws.Range('B35').Value = 'dogwifhat'
$ws.Range('C35').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.77'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '507.59'
$ws.Range('E36').Value = '  -3.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.147'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.13'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.28'
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.410'
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.16'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0837'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.26'
$ws.Range('E45').Value = '  +47.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.90'
$ws.Range('E46').Value = '  -1.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '150.72'
$ws.Range('E47').Value = '  +1.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.682'
$ws.Range('E48').Value = '  +6.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.00'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.35'
$ws.Range('E50').Value = '  +1.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.45'
$ws.Range('E51').Value = '  +1.14%  '
